# Insert a new data row into the Zanahoria (carrot) price sheet.
#
# The source data is a weekly price series for "Feria Lagunitas de Puerto
# Montt" (Zanahoria). A new week's record (date serial 44518) is inserted
# right after the existing row for date serial 44469 (row 156), pushing
# every following record down by one row. The sheet dimension grows from
# A1:R246 to A1:R247.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 156; this shifts rows 156-246
# down to 157-247 and (per Excel's normal "insert row" behavior) copies
# formatting - including the date number format on column D - from the
# row above.
$ws.Rows(156).Insert()

# Populate the newly-inserted row 156 with the new record's data. Columns
# A, B, C, E, F, G, H, I, N, Q, R are constant across every row in this
# subset, so they are simply repeated here.
$ws.Range("A156").Value = 4
$ws.Range("B156").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C156").Value = "Los Lagos"
$ws.Range("D156").Value = 44518
$ws.Range("D156").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E156").Value = 10
$ws.Range("F156").Value = 100114013
$ws.Range("G156").Value = "Zanahoria"
$ws.Range("H156").Value = "Sin especificar"
$ws.Range("I156").Value = "Primera"
$ws.Range("J156").Value = 150
$ws.Range("K156").Value = 12000
$ws.Range("L156").Value = 12000
$ws.Range("M156").Value = 12000
$ws.Range("N156").Value = "$/saco 20 kilos"
$ws.Range("O156").Value = "Región de Ñuble"
$ws.Range("P156").Value = 600
$ws.Range("Q156").Value = 20
$ws.Range("R156").Value = "Hortaliza"
